$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the NumOrden text values in column B (keep them as text with the
# trailing space, same as the originals, via the leading apostrophe so the
# quote-prefix cell style / number format is preserved).
$ws.Range("B5").Value = "'0420172010228 "
$ws.Range("B6").Value = "'1120170200969 "
$ws.Range("B7").Value = "'1220170301466 "

# The Importe values in column D for these rows are no longer populated.
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D7").ClearContents()

# Reflect the new active selection used while reviewing the edited rows.
$ws.Range("B5:B7").Select()
